$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner in A1 ---
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 21 de Abril de 2020 a las 12:52"

# --- España (row 5): F5 (Muertes hoy) updated ---
$ws.Cells.Item(5, 6).Value2 = 7705

# --- Iran (row 11): full row of stats updated ---
$ws.Cells.Item(11, 2).Value2 = 84802
$ws.Cells.Item(11, 3).Value2 = 1297
$ws.Cells.Item(11, 4).Value2 = 60965
$ws.Cells.Item(11, 5).Value2 = 18540
$ws.Cells.Item(11, 6).Value2 = 3357
$ws.Cells.Item(11, 7).Value2 = 88
$ws.Cells.Item(11, 8).Value2 = 5297

# --- Noruega (row 38): several stats updated ---
$ws.Cells.Item(38, 5).Value2 = 6942
$ws.Cells.Item(38, 6).Value2 = 56
$ws.Cells.Item(38, 7).Value2 = 1
$ws.Cells.Item(38, 8).Value2 = 182

# --- Countries around Filipinas/Bielorrusia/Ucrania/Catar/Malasia (rows 44-46) ---
# Catar's case count rose above Bielorrusia/Ucrania, so it moves up two rows in the
# (descending, sorted-by-total-cases) table while Bielorrusia/Ucrania shift down one row.
$ws.Cells.Item(44, 1).Value2 = "Catar"
$ws.Cells.Item(44, 2).Value2 = 6533
$ws.Cells.Item(44, 3).Value2 = 518
$ws.Cells.Item(44, 4).Value2 = 614
$ws.Cells.Item(44, 5).Value2 = 5910
$ws.Cells.Item(44, 6).Value2 = 37
$ws.Cells.Item(44, 7).Value2 = 0
$ws.Cells.Item(44, 8).Value2 = 9

$ws.Cells.Item(45, 1).Value2 = "Bielorrusia"
$ws.Cells.Item(45, 2).Value2 = 6264
$ws.Cells.Item(45, 3).Value2 = 0
$ws.Cells.Item(45, 4).Value2 = 514
$ws.Cells.Item(45, 5).Value2 = 5699
$ws.Cells.Item(45, 6).Value2 = 92
$ws.Cells.Item(45, 7).Value2 = 0
$ws.Cells.Item(45, 8).Value2 = 51

$ws.Cells.Item(46, 1).Value2 = "Ucrania"
$ws.Cells.Item(46, 2).Value2 = 6125
$ws.Cells.Item(46, 3).Value2 = 415
$ws.Cells.Item(46, 4).Value2 = 367
$ws.Cells.Item(46, 5).Value2 = 5597
$ws.Cells.Item(46, 6).Value2 = 45
$ws.Cells.Item(46, 7).Value2 = 10
$ws.Cells.Item(46, 8).Value2 = 161

# --- Republica de Macedonia (row 78): several stats updated ---
$ws.Cells.Item(78, 2).Value2 = 1231
$ws.Cells.Item(78, 3).Value2 = 6
$ws.Cells.Item(78, 4).Value2 = 224
$ws.Cells.Item(78, 5).Value2 = 952
$ws.Cells.Item(78, 7).Value2 = 1
$ws.Cells.Item(78, 8).Value2 = 55

# --- Madagascar (row 136): D and E updated ---
$ws.Cells.Item(136, 4).Value2 = 44
$ws.Cells.Item(136, 5).Value2 = 77

# --- Nepal (row 174): B, C and E updated ---
$ws.Cells.Item(174, 2).Value2 = 32
$ws.Cells.Item(174, 3).Value2 = 1
$ws.Cells.Item(174, 5).Value2 = 28

# --- San Vicente y las Granadinas (row 193): D and E updated ---
$ws.Cells.Item(193, 4).Value2 = 2
$ws.Cells.Item(193, 5).Value2 = 10
